$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (previously held the "Tallticka" record, now holds the "Tretåig hackspett" record)
$ws.Range("A2").Value = 80984653
$ws.Range("B2").Value = 56395
$ws.Range("E2").Value = 100109
$ws.Range("F2").Value = "Tretåig hackspett"
$ws.Range("G2").Value = "Picoides tridactylus"
$ws.Range("H2").Value = "(Linnaeus, 1758)"
$ws.Range("Q2").Value = 802412.2482731647
$ws.Range("R2").Value = 7420062.223471683
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2019-06-26"
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2019-06-26"

# Row 3 (keeps "Skinnlav" but with different Id/coordinates)
$ws.Range("A3").Value = 80984665
$ws.Range("Q3").Value = 802408.9692803487
$ws.Range("R3").Value = 7420170.958352051

# Row 4 (previously "Tretåig hackspett", now "Garnlav")
$ws.Range("A4").Value = 80984695
$ws.Range("B4").Value = 77506
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = "Garnlav"
$ws.Range("G4").Value = "Alectoria sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."

# Row 5 (previously "Skinnlav", now "Jungfru marie nycklar"; Auktor becomes blank)
$ws.Range("A5").Value = 80984688
$ws.Range("B5").Value = 96254
$ws.Range("E5").Value = 223597
$ws.Range("F5").Value = "Jungfru marie nycklar"
$ws.Range("G5").Value = "Dactylorhiza maculata subsp. maculata"
$ws.Range("H5").Value = ""
$ws.Range("Q5").Value = 802414.0848113028
$ws.Range("R5").Value = 7420060.03789904

# Row 6 (previously "Garnlav", now "Tallticka")
$ws.Range("A6").Value = 80984727
$ws.Range("B6").Value = 89412
$ws.Range("E6").Value = 5442
$ws.Range("F6").Value = "Tallticka"
$ws.Range("G6").Value = "Porodaedalea pini"
$ws.Range("H6").Value = "(Brot.) Murrill"
$ws.Range("Q6").Value = 802488.0892767747
$ws.Range("R6").Value = 7419795.075097569
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2019-06-25"
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2019-06-25"

# Row 8 (previously "Jungfru marie nycklar", now "Skinnlav")
$ws.Range("A8").Value = 80984668
$ws.Range("B8").Value = 78503
$ws.Range("E8").Value = 6456
$ws.Range("F8").Value = "Skinnlav"
$ws.Range("G8").Value = "Leptogium saturninum"
$ws.Range("H8").Value = "(Dicks.) Nyl."
$ws.Range("Q8").Value = 802997.0783882558
$ws.Range("R8").Value = 7419415.035970301
